$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "‘민족의 미래를 사겠다는 장사꾼’"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/merchant-for-korean-future/#utm_source=rss&utm_medium=rss&utm_campaign=merchant-for-korean-future"

$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"

$ws.Range("D51").Value = "[pandas] 파이썬 판다스로 엑셀 파일을 읽고 쓰려면 openpyxl도 추가로 설치해야 함"
$ws.Range("E51").Value = "https://bskyvision.com/entry/pandas-%ED%8C%8C%EC%9D%B4%EC%8D%AC-%ED%8C%90%EB%8B%A4%EC%8A%A4%EB%A1%9C-%EC%97%91%EC%85%80-%ED%8C%8C%EC%9D%BC%EC%9D%84-%EC%9D%BD%EA%B3%A0-%EC%93%B0%EB%A0%A4%EB%A9%B4-openpyxl%EB%8F%84-%EC%B6%94%EA%B0%80%EB%A1%9C-%EC%84%A4%EC%B9%98%ED%95%B4%EC%95%BC-%ED%95%A8"
